# Insert a new data row at row 527 (shifts existing rows 527-600 down to 528-601)
# and populate it with the new "Choclo" price-report record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(527).Insert()

$ws.Cells.Item(527, 1).Value  = 9
$ws.Cells.Item(527, 2).Value  = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(527, 3).Value  = 'Metropolitana'
$ws.Cells.Item(527, 4).Value  = 44816
$ws.Cells.Item(527, 5).Value  = 13
$ws.Cells.Item(527, 6).Value  = 100112024
$ws.Cells.Item(527, 7).Value  = 'Choclo'
$ws.Cells.Item(527, 8).Value  = 'Dulce o Americano'
$ws.Cells.Item(527, 9).Value  = 'Primera'
$ws.Cells.Item(527, 10).Value = 110
$ws.Cells.Item(527, 11).Value = 19000
$ws.Cells.Item(527, 12).Value = 20000
$ws.Cells.Item(527, 13).Value = 19455
$ws.Cells.Item(527, 14).Value = '$/malla 70 unidades'
$ws.Cells.Item(527, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(527, 16).Value = 278
$ws.Cells.Item(527, 17).Value = 70
$ws.Cells.Item(527, 18).Value = 'Hortaliza'
